# Auto-generated: refresh cached Universalis market-price columns (H:N) per row
# across all 8 item-category sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values come from a scheduled data-refresh run; no formulas are involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 60.57143
$ws.Range("I11").Value = 60.57143
$ws.Range("K11").Value = 60.57143
$ws.Range("M11").Value = 79.42857000000001
$ws.Range("H17").Value = 1023322.2
$ws.Range("J17").Value = 1614349.5
$ws.Range("L17").Value = 4843048.5
$ws.Range("N17").Value = -4843384.5
$ws.Range("H137").Value = 1157.1082
$ws.Range("I137").Value = 975.7959
$ws.Range("J137").Value = 1512.48
$ws.Range("K137").Value = 2927.3877
$ws.Range("L137").Value = 4537.440000000001
$ws.Range("M137").Value = -377.3876999999998
$ws.Range("N137").Value = -9637.440000000001
$ws.Range("H138").Value = 2240
$ws.Range("I138").Value = 1481.7632
$ws.Range("J138").Value = 2880.2888
$ws.Range("K138").Value = 4445.2896
$ws.Range("L138").Value = 8640.866399999999
$ws.Range("M138").Value = 694.7103999999999
$ws.Range("N138").Value = -18920.8664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H32").Value = 914234.1
$ws.Range("I32").Value = 1048232.1
$ws.Range("J32").Value = 20914.223
$ws.Range("K32").Value = 1048232.1
$ws.Range("L32").Value = 20914.223
$ws.Range("M32").Value = -1047945.1
$ws.Range("N32").Value = -21488.223
$ws.Range("H45").Value = 2504.2727
$ws.Range("I45").Value = 2530.875
$ws.Range("J45").Value = 2433.3333
$ws.Range("K45").Value = 2530.875
$ws.Range("L45").Value = 2433.3333
$ws.Range("M45").Value = -2153.875
$ws.Range("N45").Value = -3187.3333
$ws.Range("H61").Value = 6291182.5
$ws.Range("I61").Value = 7753399
$ws.Range("J61").Value = 3651.3
$ws.Range("K61").Value = 7753399
$ws.Range("L61").Value = 3651.3
$ws.Range("M61").Value = -7753187
$ws.Range("N61").Value = -4075.3
$ws.Range("H74").Value = 703.7941
$ws.Range("I74").Value = 703.7941
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 703.7941
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 170.2059
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 703.7941
$ws.Range("I77").Value = 703.7941
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 3518.9705
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 849.0295000000001
$ws.Range("N77").ClearContents()
$ws.Range("H80").Value = 19297.334
$ws.Range("J80").Value = 19297.334
$ws.Range("L80").Value = 19297.334
$ws.Range("N80").Value = -21293.334
$ws.Range("H83").Value = 19297.334
$ws.Range("J83").Value = 19297.334
$ws.Range("L83").Value = 57892.00199999999
$ws.Range("N83").Value = -67876.00199999999
$ws.Range("H133").Value = 30030.5
$ws.Range("J133").Value = 30030.5
$ws.Range("L133").Value = 30030.5
$ws.Range("N133").Value = -35090.5
$ws.Range("H135").Value = 38163.25
$ws.Range("J135").Value = 38163.25
$ws.Range("L135").Value = 38163.25
$ws.Range("N135").Value = -48303.25
$ws.Range("H136").Value = 6291182.5
$ws.Range("I136").Value = 7753399
$ws.Range("J136").Value = 3651.3
$ws.Range("K136").Value = 23260197
$ws.Range("L136").Value = 10953.9
$ws.Range("M136").Value = -23257647
$ws.Range("N136").Value = -16053.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1802.8846
$ws.Range("I20").Value = 1439.8667
$ws.Range("J20").Value = 2297.9092
$ws.Range("K20").Value = 1439.8667
$ws.Range("L20").Value = 2297.9092
$ws.Range("M20").Value = -1192.8667
$ws.Range("N20").Value = -2791.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4313.684
$ws.Range("I31").Value = 1358.7
$ws.Range("J31").Value = 7597
$ws.Range("K31").Value = 1358.7
$ws.Range("L31").Value = 7597
$ws.Range("M31").Value = -1063.7
$ws.Range("N31").Value = -8187
$ws.Range("H34").Value = 4313.684
$ws.Range("I34").Value = 1358.7
$ws.Range("J34").Value = 7597
$ws.Range("K34").Value = 1358.7
$ws.Range("L34").Value = 7597
$ws.Range("M34").Value = -1156.7
$ws.Range("N34").Value = -8001
$ws.Range("H98").Value = 39999
$ws.Range("J98").Value = 39999
$ws.Range("L98").Value = 39999
$ws.Range("N98").Value = -44491
$ws.Range("H132").Value = 5954439.5
$ws.Range("I132").Value = 1837.2354
$ws.Range("K132").Value = 5511.706200000001
$ws.Range("M132").Value = -2981.706200000001
$ws.Range("H134").Value = 3771
$ws.Range("I134").Value = 3618.9744
$ws.Range("J134").Value = 4759.1665
$ws.Range("K134").Value = 10856.9232
$ws.Range("L134").Value = 14277.4995
$ws.Range("M134").Value = -8321.923200000001
$ws.Range("N134").Value = -19347.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 43478388
$ws.Range("I23").Value = 122
$ws.Range("J23").Value = 55555684
$ws.Range("K23").Value = 366
$ws.Range("L23").Value = 166667052
$ws.Range("M23").Value = -131
$ws.Range("N23").Value = -166667522
$ws.Range("H34").Value = 10000450
$ws.Range("I34").Value = 170.6
$ws.Range("J34").Value = 11111592
$ws.Range("K34").Value = 511.8
$ws.Range("L34").Value = 33334776
$ws.Range("M34").Value = -427.8
$ws.Range("N34").Value = -33334944
$ws.Range("H69").Value = 1417.6666
$ws.Range("I69").Value = 603
$ws.Range("J69").Value = 1825
$ws.Range("K69").Value = 1809
$ws.Range("L69").Value = 5475
$ws.Range("M69").Value = -998
$ws.Range("N69").Value = -7097
$ws.Range("H72").Value = 1417.6666
$ws.Range("I72").Value = 603
$ws.Range("J72").Value = 1825
$ws.Range("K72").Value = 5427
$ws.Range("L72").Value = 16425
$ws.Range("M72").Value = -1371
$ws.Range("N72").Value = -24537
$ws.Range("H107").Value = 15384837
$ws.Range("I107").Value = 286.90323
$ws.Range("J107").Value = 29411928
$ws.Range("K107").Value = 860.70969
$ws.Range("L107").Value = 88235784
$ws.Range("M107").Value = 1059.29031
$ws.Range("N107").Value = -88239624
$ws.Range("H137").Value = 4908095
$ws.Range("J137").Value = 4055.2104
$ws.Range("L137").Value = 12165.6312
$ws.Range("N137").Value = -22365.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -360
$ws.Range("N12").ClearContents()
$ws.Range("H80").Value = 15921307
$ws.Range("I80").Value = 16969294
$ws.Range("J80").Value = 201500
$ws.Range("K80").Value = 16969294
$ws.Range("L80").Value = 201500
$ws.Range("M80").Value = -16968296
$ws.Range("N80").Value = -203496
$ws.Range("H83").Value = 15921307
$ws.Range("I83").Value = 16969294
$ws.Range("J83").Value = 201500
$ws.Range("K83").Value = 84846470
$ws.Range("L83").Value = 1007500
$ws.Range("M83").Value = -84841478
$ws.Range("N83").Value = -1017484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 26693.666
$ws.Range("J101").Value = 26693.666
$ws.Range("L101").Value = 26693.666
$ws.Range("N101").Value = -33183.666
$ws.Range("H132").Value = 2480.0396
$ws.Range("I132").Value = 2234.3447
$ws.Range("J132").Value = 3271.7222
$ws.Range("K132").Value = 6703.034100000001
$ws.Range("L132").Value = 9815.1666
$ws.Range("M132").Value = -4173.034100000001
$ws.Range("N132").Value = -14875.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 200
$ws.Range("J8").Value = 300
$ws.Range("L8").Value = 300
$ws.Range("N8").Value = -580
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H132").Value = 4228286.5
$ws.Range("I132").Value = 1162.5869
$ws.Range("J132").Value = 12682535
$ws.Range("K132").Value = 3487.7607
$ws.Range("L132").Value = 38047605
$ws.Range("M132").Value = -957.7606999999998
$ws.Range("N132").Value = -38052665
$ws.Range("H136").Value = 936.8200000000001
$ws.Range("I136").Value = 879.0633
$ws.Range("J136").Value = 1154.0952
$ws.Range("K136").Value = 2637.1899
$ws.Range("L136").Value = 3462.2856
$ws.Range("M136").Value = -87.18990000000031
$ws.Range("N136").Value = -8562.285599999999
